$wb = $excel.ActiveWorkbook

# Sheet2 has the duplicate "b" header in D1 - fix it to "c"
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("D1").Value = "c"

# Update the selection on Sheet2 to F7 (also marks it the active sheet/view)
$ws2.Activate()
$ws2.Range("F7").Select()
